$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - player / withdraw record (ID, name, withdraw date, source account)
$ws.Range("A2").Value = "angelboyx"
$ws.Range("B2").Value = "ASEP SUPRIYADI"
$ws.Range("C2").Value = "'2022-04-24"
$ws.Range("F2").Value = 12345665

# Row 3 - player / withdraw record (ID, name, withdraw date, amount, source account)
$ws.Range("A3").Value = "player11"
$ws.Range("B3").Value = "Bejo Kuncoro"
$ws.Range("C3").Value = "'2022-04-24"
$ws.Range("D3").Value = 15000000
$ws.Range("F3").Value = 12345665

# Selection moves from F15 to D4
$null = $ws.Range("D4").Select()

# Row 1 grows slightly taller (bold header re-measured) once the content changed.
$ws.Rows.Item(1).RowHeight = 15.6

# Columns re-measure to "best fit" the new, longer text (names/dates/numbers).
# target widths (stored OOXML chars): 9.33203125 12.77734375 18.21875 17.77734375 12.6640625 23
$ws.Columns.Item(1).ColumnWidth = 8.5
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(3).ColumnWidth = 17.333333333333332
$ws.Columns.Item(4).ColumnWidth = 17
$ws.Columns.Item(5).ColumnWidth = 11.833333333333334
$ws.Columns.Item(6).ColumnWidth = 22.166666666666668
